# Rebuild the "plot_infos_pcr_msc" sheet:
#  - header columns C:F are reordered (R², RMSE, Offset, Slope)
#  - every attribute (SST, PH, AT, FIRMEZA (N), UBS (%)) gains a third
#    "Validação" row in addition to "Referência"/"Predição", and the
#    metric values are refreshed, expanding the table from A1:F11 to A1:F16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (reordered): C1=R², D1=RMSE, E1=Offset, F1=Slope
$ws.Range("C1").Value = "R²"
$ws.Range("D1").Value = "RMSE"
$ws.Range("E1").Value = "Offset"
$ws.Range("F1").Value = "Slope"

# Data rows 2-16
$ws.Range("A2").Value = "SST"
$ws.Range("B2").Value = "Referência"
$ws.Range("C2").Value = 0.7066831166592163
$ws.Range("D2").Value = 1.438670180008751
$ws.Range("E2").Value = 4.105095489589983
$ws.Range("F2").Value = 0.7066831166592163

$ws.Range("A3").Value = "SST"
$ws.Range("B3").Value = "Predição"
$ws.Range("C3").Value = 0.6804698690085335
$ws.Range("D3").Value = 1.501580499895086
$ws.Range("E3").Value = 4.248792685759586
$ws.Range("F3").Value = 0.6971472239507518

$ws.Range("A4").Value = "SST"
$ws.Range("B4").Value = "Validação"
$ws.Range("C4").Value = 0.582392765366923
$ws.Range("D4").Value = 1.359874398235301
$ws.Range("E4").Value = 1.764861058949153
$ws.Range("F4").Value = 0.8890481639202232

$ws.Range("A5").Value = "PH"
$ws.Range("B5").Value = "Referência"
$ws.Range("C5").Value = 0.4246018486687569
$ws.Range("D5").Value = 0.2339524422963989
$ws.Range("E5").Value = 1.889906736010495
$ws.Range("F5").Value = 0.4246018486687571

$ws.Range("A6").Value = "PH"
$ws.Range("B6").Value = "Predição"
$ws.Range("C6").Value = 0.3175427648022607
$ws.Range("D6").Value = 0.2547892341226295
$ws.Range("E6").Value = 2.04525999503011
$ws.Range("F6").Value = 0.3770829326099566

$ws.Range("A7").Value = "PH"
$ws.Range("B7").Value = "Validação"
$ws.Range("C7").Value = 0.4193890325460984
$ws.Range("D7").Value = 0.2008072895261398
$ws.Range("E7").Value = 0.6091880822025666
$ws.Range("F7").Value = 0.8234431080378893

$ws.Range("A8").Value = "AT"
$ws.Range("B8").Value = "Referência"
$ws.Range("C8").Value = 0.4907501243441826
$ws.Range("D8").Value = 0.4184808130573737
$ws.Range("E8").Value = 0.5717741203891936
$ws.Range("F8").Value = 0.4907501243441822

$ws.Range("A9").Value = "AT"
$ws.Range("B9").Value = "Predição"
$ws.Range("C9").Value = 0.4378967930685279
$ws.Range("D9").Value = 0.4396611752888547
$ws.Range("E9").Value = 0.595713707421827
$ws.Range("F9").Value = 0.468931609688125

$ws.Range("A10").Value = "AT"
$ws.Range("B10").Value = "Validação"
$ws.Range("C10").Value = 0.4715904908622609
$ws.Range("D10").Value = 0.3179403569792187
$ws.Range("E10").Value = -0.01515743881380183
$ws.Range("F10").Value = 0.951333887657626

$ws.Range("A11").Value = "FIRMEZA (N)"
$ws.Range("B11").Value = "Referência"
$ws.Range("C11").Value = 0.468851957774996
$ws.Range("D11").Value = 74.8986845060683
$ws.Range("E11").Value = 272.4843057039558
$ws.Range("F11").Value = 0.4688519577749965

$ws.Range("A12").Value = "FIRMEZA (N)"
$ws.Range("B12").Value = "Predição"
$ws.Range("C12").Value = 0.4280222426197769
$ws.Range("D12").Value = 77.72414758672673
$ws.Range("E12").Value = 282.0857255236262
$ws.Range("F12").Value = 0.4496581277236957

$ws.Range("A13").Value = "FIRMEZA (N)"
$ws.Range("B13").Value = "Validação"
$ws.Range("C13").Value = 0.4698552868879904
$ws.Range("D13").Value = 56.25050462134656
$ws.Range("E13").Value = 3.660049469787396
$ws.Range("F13").Value = 0.9676263846128885

$ws.Range("A14").Value = "UBS (%)"
$ws.Range("B14").Value = "Referência"
$ws.Range("C14").Value = 0.6324113483372676
$ws.Range("D14").Value = 1.950333729245592
$ws.Range("E14").Value = 5.638455527228118
$ws.Range("F14").Value = 0.6324113483372673

$ws.Range("A15").Value = "UBS (%)"
$ws.Range("B15").Value = "Predição"
$ws.Range("C15").Value = 0.5946897203807606
$ws.Range("D15").Value = 2.047961053186621
$ws.Range("E15").Value = 5.878397254636665
$ws.Range("F15").Value = 0.6172330804647264

$ws.Range("A16").Value = "UBS (%)"
$ws.Range("B16").Value = "Validação"
$ws.Range("C16").Value = 0.6326836387267732
$ws.Range("D16").Value = 1.313405097363099
$ws.Range("E16").Value = 2.841349805743501
$ws.Range("F16").Value = 0.8253233738178234
